$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2300.125
$ws.Range("J51").Value = 2400
$ws.Range("L51").Value = 2400
$ws.Range("N51").Value = -3368
$ws.Range("H69").Value = 3822.8572
$ws.Range("I69").Value = 5850
$ws.Range("J69").Value = 3012
$ws.Range("K69").Value = 17550
$ws.Range("L69").Value = 9036
$ws.Range("M69").Value = -16676
$ws.Range("N69").Value = -10784
$ws.Range("H72").Value = 3822.8572
$ws.Range("I72").Value = 5850
$ws.Range("J72").Value = 3012
$ws.Range("K72").Value = 52650
$ws.Range("L72").Value = 27108
$ws.Range("M72").Value = -48282
$ws.Range("N72").Value = -35844
$ws.Range("H106").Value = 2284.6155
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H138").Value = 3042.7
$ws.Range("I138").Value = 1930
$ws.Range("J138").Value = 3578.4443
$ws.Range("K138").Value = 5790
$ws.Range("L138").Value = 10735.3329
$ws.Range("M138").Value = -650
$ws.Range("N138").Value = -21015.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H24").Value = 32500
$ws.Range("J24").Value = 32500
$ws.Range("L24").Value = 32500
$ws.Range("N24").Value = -33248
$ws.Range("H100").Value = 32500
$ws.Range("J100").Value = 32500
$ws.Range("L100").Value = 32500
$ws.Range("N100").Value = -34664
$ws.Range("H102").Value = 3707293.5
$ws.Range("I102").Value = 6175889.5
$ws.Range("K102").Value = 6175889.5
$ws.Range("M102").Value = -6174267.5
$ws.Range("H122").Value = 7354273.5
$ws.Range("I122").Value = 1286.5834
$ws.Range("J122").Value = 25001442
$ws.Range("K122").Value = 3859.7502
$ws.Range("L122").Value = 75004326
$ws.Range("M122").Value = -1409.7502
$ws.Range("N122").Value = -75009226

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1751
$ws.Range("I20").Value = 1974.2
$ws.Range("J20").Value = 1332.5
$ws.Range("K20").Value = 1974.2
$ws.Range("L20").Value = 1332.5
$ws.Range("M20").Value = -1727.2
$ws.Range("N20").Value = -1826.5
$ws.Range("H134").Value = 57175.5
$ws.Range("I134").Value = 1714.4117
$ws.Range("K134").Value = 5143.2351
$ws.Range("M134").Value = -2608.2351

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 566666.75
$ws.Range("J4").Value = 566666.75
$ws.Range("L4").Value = 566666.75
$ws.Range("N4").Value = -566890.75
$ws.Range("H86").Value = 3577.8
$ws.Range("I86").Value = 4250
$ws.Range("J86").Value = 3129.6667
$ws.Range("K86").Value = 4250
$ws.Range("L86").Value = 3129.6667
$ws.Range("M86").Value = -3127
$ws.Range("N86").Value = -5375.6667
$ws.Range("H89").Value = 3577.8
$ws.Range("I89").Value = 4250
$ws.Range("J89").Value = 3129.6667
$ws.Range("K89").Value = 21250
$ws.Range("L89").Value = 15648.3335
$ws.Range("M89").Value = -15634
$ws.Range("N89").Value = -26880.3335
$ws.Range("H94").Value = 1000
$ws.Range("J94").Value = 1000
$ws.Range("L94").Value = 1000
$ws.Range("N94").Value = -1902
$ws.Range("H105").Value = 682.8823
$ws.Range("I105").Value = 663.0625
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 663.0625
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 1083.9375
$ws.Range("N105").Value = -4494
$ws.Range("H107").Value = 534.1
$ws.Range("I107").Value = 330.85715
$ws.Range("J107").Value = 1008.3333
$ws.Range("K107").Value = 330.85715
$ws.Range("L107").Value = 1008.3333
$ws.Range("M107").Value = 1589.14285
$ws.Range("N107").Value = -4848.3333
$ws.Range("H122").Value = 5808.4688
$ws.Range("I122").Value = 2419.875
$ws.Range("K122").Value = 7259.625
$ws.Range("M122").Value = -4809.625
$ws.Range("H134").Value = 2832.75
$ws.Range("I134").Value = 2404.238
$ws.Range("J134").Value = 4118.2856
$ws.Range("K134").Value = 7212.714
$ws.Range("L134").Value = 12354.8568
$ws.Range("M134").Value = -4677.714
$ws.Range("N134").Value = -17424.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 10423246
$ws.Range("I5").Value = 366.15384
$ws.Range("J5").Value = 55589060
$ws.Range("K5").Value = 1098.46152
$ws.Range("L5").Value = 166767180
$ws.Range("M5").Value = -986.4615200000001
$ws.Range("N5").Value = -166767404
$ws.Range("H12").Value = 66666944
$ws.Range("I12").Value = 200000110
$ws.Range("J12").Value = 359.6
$ws.Range("K12").Value = 600000330
$ws.Range("L12").Value = 1078.8
$ws.Range("M12").Value = -600000157
$ws.Range("N12").Value = -1424.8
$ws.Range("H68").Value = 2571.4062
$ws.Range("I68").Value = 924.8387
$ws.Range("J68").Value = 4118.1816
$ws.Range("K68").Value = 2774.5161
$ws.Range("L68").Value = 12354.5448
$ws.Range("M68").Value = -1963.5161
$ws.Range("N68").Value = -13976.5448
$ws.Range("H71").Value = 2571.4062
$ws.Range("I71").Value = 924.8387
$ws.Range("J71").Value = 4118.1816
$ws.Range("K71").Value = 8323.5483
$ws.Range("L71").Value = 37063.6344
$ws.Range("M71").Value = -4267.5483
$ws.Range("N71").Value = -45175.6344
$ws.Range("H92").Value = 333.33334
$ws.Range("I92").Value = 250
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 750
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 498
$ws.Range("N92").Value = -3996
$ws.Range("H107").Value = 765.37933
$ws.Range("I107").Value = 290.1905
$ws.Range("K107").Value = 870.5715
$ws.Range("M107").Value = 1049.4285
$ws.Range("H113").Value = 429.91
$ws.Range("I113").Value = 485.68
$ws.Range("K113").Value = 1457.04
$ws.Range("M113").Value = 712.96
$ws.Range("H122").Value = 808.3125
$ws.Range("I122").Value = 237.11111
$ws.Range("J122").Value = 1542.7142
$ws.Range("K122").Value = 2133.99999
$ws.Range("L122").Value = 13884.4278
$ws.Range("M122").Value = 316.0000100000002
$ws.Range("N122").Value = -18784.4278
$ws.Range("H135").Value = 10423246
$ws.Range("I135").Value = 366.15384
$ws.Range("J135").Value = 55589060
$ws.Range("K135").Value = 3295.38456
$ws.Range("L135").Value = 500301540
$ws.Range("M135").Value = -760.38456
$ws.Range("N135").Value = -500306610

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1559.4546
$ws.Range("I97").Value = 1559.4546
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1559.4546
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1063.4546
$ws.Range("N97").ClearContents()
$ws.Range("H113").Value = 2049.0435
$ws.Range("I113").Value = 2023.9231
$ws.Range("J113").Value = 2081.7
$ws.Range("K113").Value = 2023.9231
$ws.Range("L113").Value = 2081.7
$ws.Range("M113").Value = 146.0769
$ws.Range("N113").Value = -6421.7
$ws.Range("H126").Value = 2795.238
$ws.Range("I126").Value = 1911.1111
$ws.Range("J126").Value = 3458.3333
$ws.Range("K126").Value = 5733.3333
$ws.Range("L126").Value = 10374.9999
$ws.Range("M126").Value = -3263.3333
$ws.Range("N126").Value = -15314.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 63550
$ws.Range("J92").Value = 63550
$ws.Range("L92").Value = 63550
$ws.Range("N92").Value = -68542
$ws.Range("H107").Value = 1944
$ws.Range("I107").Value = 1171.375
$ws.Range("K107").Value = 3514.125
$ws.Range("M107").Value = -1594.125
